# The document's SharePoint-linked custom XML parts get renumbered:
# customXml/item2.xml (p:properties / documentManagement) and
# customXml/item4.xml (ct:contentTypeSchema) swap slots, while
# customXml/item3.xml (the bibliography b:Sources part) keeps its slot.
# Pull each part out by its root namespace, delete the pair, and re-add
# them in the swapped order so the package is re-serialized with the new
# item numbering (contentTypeSchema first, documentManagement properties
# last).

$d = $word.ActiveDocument
$parts = $d.CustomXMLParts

$propsNS  = "http://schemas.microsoft.com/office/2006/metadata/properties"
$schemaNS = "http://schemas.microsoft.com/office/2006/metadata/contentType"

$propsPart  = $null
$schemaPart = $null
$propsXML   = $null
$schemaXML  = $null

for ($i = 1; $i -le $parts.Count; $i++) {
    $p = $parts.Item($i)
    if ($p.NamespaceURI -eq $propsNS) {
        $propsPart = $p
        $propsXML = $p.XML
    } elseif ($p.NamespaceURI -eq $schemaNS) {
        $schemaPart = $p
        $schemaXML = $p.XML
    }
}

if ($schemaPart -ne $null) { $schemaPart.Delete() }
if ($propsPart -ne $null) { $propsPart.Delete() }

if ($schemaXML -ne $null) {
    # content-type schema now occupies the earlier item slot
    [void]$parts.Add($schemaXML)
}
if ($propsXML -ne $null) {
    # document-management properties now occupy the later item slot
    [void]$parts.Add($propsXML)
}

$d.Saved = $false
$d.Save()
